# The admin page's prodText/productPrice/disCountPrice columns were
# producing bad data (page conflict), so this page (used by the track
# order test) now just reports "null" for those three columns on every
# product row. imgsrc / prodLink (columns A and B) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 3).Value = "null"
    $ws.Cells.Item($r, 4).Value = "null"
    $ws.Cells.Item($r, 5).Value = "null"
}
